# AutomationControlSheet.xlsx - test case name updated as per documentation
#
# 1) AppControl sheet: B25 mailing-list cell is replaced with a single
#    e-mail address and turned into a (mailto:) hyperlink; the row is
#    shrunk to fit the now much shorter wrapped text, and the sheet's
#    selection moves on to A26.
# 2) smoke sheet: the seven automated test-case names in A17:A23 get a
#    "ZestIOT_" prefix (per the updated documentation/naming convention),
#    and rows 20/21 grow a bit to fit the now-longer wrapped text.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# AppControl
# ---------------------------------------------------------------
$wsApp = $wb.Worksheets.Item("AppControl")

$wsApp.Range("B25").Value = "stiyyagura@enhops.com"
$wsApp.Hyperlinks.Add($wsApp.Range("B25"), "mailto:stiyyagura@enhops.com") | Out-Null
$wsApp.Rows.Item(25).RowHeight = 45

$wsApp.Range("A26").Select()

# ---------------------------------------------------------------
# smoke
# ---------------------------------------------------------------
$wsSmoke = $wb.Worksheets.Item("smoke")

# Re-activate "smoke" so it stays the workbook's active/selected tab
# (only AppControl's own selected cell is moving on to A26).
$wsSmoke.Activate()

$wsSmoke.Range("A17").Value = "ZestIOT_AV_2268_Validate_Accuracy_of_COBT_For_DIALCelebi_User"
$wsSmoke.Range("A18").Value = "ZestIOT_AV_2268_Validate_Accuracy_of_COBT_For_GMR_HYD_AISATS_User"
$wsSmoke.Range("A19").Value = "ZestIOT_AV_2268_Validate_Accuracy_of_COBT_For_GMR_HYD_SG_User"
$wsSmoke.Range("A20").Value = "ZestIOT_AV_2293_Identify_coverage_of_Flight_Sensor_and_Validate_timestamps_of_Arrival_Aircrafts"
$wsSmoke.Range("A21").Value = "ZestIOT_AV_2294_Identify_coverage_of_Flight_Sensor_and_Validate_timestamps_of_Departure_Aircrafts"
$wsSmoke.Range("A22").Value = "ZestIOT_AV_2307_Validate_LANDING_ONBLOCK_OFFBLOCK_AIRBORNE_timestamps_of_Arrival_and_Departure_aircrafts_Any_Data_source"
$wsSmoke.Range("A23").Value = "ZestIOT_AV_2304_Identify_the_coverage_of_Boarding_activities_and_validate_timestamps"

$wsSmoke.Rows.Item(20).RowHeight = 30
$wsSmoke.Rows.Item(21).RowHeight = 45
